$wb = $excel.ActiveWorkbook

# --- Update "queries" sheet: insert new column D "linked_form_id" ---
$queries = $wb.Worksheets.Item("queries")
$queries.Columns("D").Insert()
$queries.Range("D1").Value = "linked_form_id"
$queries.Range("D2").Value = "entitlements"

# --- Add new "initial" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "initial"

$ws.Range("A1").Value = "clause"
$ws.Range("B1").Value = "type"
$ws.Range("C1").Value = "display.text"
$ws.Range("D1").Value = "comments"
$ws.Range("A2").Value = "do section survey"

$ws.Range("A1:D2").WrapText = $true
$ws.Rows("3:3").EntireRow.Select() | Out-Null

$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

Write-Output "done"
